$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect numeric-looking Price (column D) values as text so they keep their
# original formatting (trailing zeros, thousand-dot grouping, etc.) instead of
# being auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '64.326.22'
$ws.Range("E2").Value = '  -4.01%  '
$ws.Range("D3").Value = '3.161.42'
$ws.Range("E3").Value = '  -8.49%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '563.89'
$ws.Range("E5").Value = '  -4.01%  '
$ws.Range("D6").Value = '170.48'
$ws.Range("E6").Value = '  -5.23%  '
$ws.Range("D7").Value = '0.612'
$ws.Range("E7").Value = '  -3.07%  '
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("D9").Value = '3.156.91'
$ws.Range("E9").Value = '  -8.60%  '
$ws.Range("E10").Value = '  -7.55%  '
$ws.Range("E11").Value = '  -6.03%  '
$ws.Range("D12").Value = '0.394'
$ws.Range("E12").Value = '  -6.10%  '
$ws.Range("D13").Value = '3.706.20'
$ws.Range("E13").Value = '  -8.58%  '
$ws.Range("E14").Value = '  +0.47%  '
$ws.Range("D15").Value = '27.21'
$ws.Range("E15").Value = '  -9.91%  '
$ws.Range("D16").Value = '64.246.73'
$ws.Range("E16").Value = '  -4.00%  '
$ws.Range("D17").Value = '0.0000162'
$ws.Range("E17").Value = '  -7.39%  '
$ws.Range("D18").Value = '3.155.59'
$ws.Range("E18").Value = '  -8.57%  '
$ws.Range("D19").Value = '5.73'
$ws.Range("E19").Value = '  -4.09%  '
$ws.Range("D20").Value = '12.94'
$ws.Range("E20").Value = '  -7.15%  '
$ws.Range("D21").Value = '354.20'
$ws.Range("E21").Value = '  -5.42%  '
$ws.Range("D22").Value = '7.21'
$ws.Range("E22").Value = '  -6.39%  '
$ws.Range("D23").Value = '1.00'
$ws.Range("E23").Value = '  +0.27%  '
$ws.Range("D24").Value = '68.07'
$ws.Range("E24").Value = '  -7.66%  '
$ws.Range("D25").Value = '0.500'
$ws.Range("E25").Value = '  -7.14%  '
$ws.Range("D26").Value = '0.0000117'
$ws.Range("E26").Value = '  -11.05%  '
$ws.Range("D27").Value = '9.61'
$ws.Range("E27").Value = '  -4.28%  '
$ws.Range("E28").Value = '  -2.57%  '
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  -0.12%  '
$ws.Range("E30").Value = '  -0.18%  '
$ws.Range("D31").Value = '1.89'
$ws.Range("E31").Value = '  -5.68%  '
$ws.Range("D32").Value = '5.39'
$ws.Range("E32").Value = '  -9.02%  '
$ws.Range("D33").Value = '21.94'
$ws.Range("E33").Value = '  -7.74%  '
$ws.Range("E34").Value = '  -6.98%  '
$ws.Range("D35").Value = '6.64'
$ws.Range("E35").Value = '  -7.11%  '
$ws.Range("E36").Value = '  -9.59%  '
$ws.Range("D37").Value = '153.83'
$ws.Range("E37").Value = '  -5.61%  '
$ws.Range("D38").Value = '0.824'
$ws.Range("E38").Value = '  -6.98%  '
$ws.Range("D39").Value = '26.15'
$ws.Range("E39").Value = '  -6.74%  '
$ws.Range("D40").Value = '1.70'
$ws.Range("E40").Value = '  -6.89%  '
$ws.Range("B41").Value = 'dogwifhat'
$ws.Range("C41").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D41").Value = '2.47'
$ws.Range("E41").Value = '  -7.15%  '
$ws.Range("B42").Value = 'Maker'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D42").Value = '2.641.71'
$ws.Range("E42").Value = '  -4.47%  '
$ws.Range("E43").Value = '  -8.08%  '
$ws.Range("B44").Value = 'OKB'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D44").Value = '39.24'
$ws.Range("E44").Value = '  -2.47%  '
$ws.Range("B45").Value = 'RenderToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D45").Value = '5.94'
$ws.Range("E45").Value = '  -7.13%  '
$ws.Range("D46").Value = '0.0651'
$ws.Range("E46").Value = '  -7.18%  '
$ws.Range("D47").Value = '23.98'
$ws.Range("E47").Value = '  -6.99%  '
$ws.Range("D48").Value = '317.83'
$ws.Range("E48").Value = '  -6.50%  '
$ws.Range("D49").Value = '0.0272'
$ws.Range("E49").Value = '  -6.11%  '
$ws.Range("E50").Value = '  -4.46%  '
$ws.Range("E51").Value = '  -0.18%  '

# Restore the default (unstyled) cell style now that the text values are locked in,
# so the saved workbook does not pick up a stray "@" number format on column D.
$ws.Range("D2:D51").Style = "Normal"

